# Update countries & provincias Spain
# Applies the covid-dashboard refresh: a handful of country rows get new
# totals, three country pairs swap places in the ranking (their row's
# country name + figures exchange), and the "last updated" timestamp moves
# from 08:58 to 10:15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 10:15"

# --- 2. Plain numeric refreshes (country stays on its row) ---------------
# row, B, C, D, E, F, G, H
$plainRows = @(
    @(6,   2464316, 4703, 1752829, 663310, 0, 33,  48177),
    @(7,    912823, 5065,  722964, 174361, 0, 114, 15498),
    @(25,   153660, 6216,   71405,  79813, 0, 16,   2442),
    @(34,    87872, 1732,   46797,  39064, 0, 19,   2011),
    @(47,    55580,   83,   50736,   4817, 0,  0,     27),
    @(71,    22743,  385,   13350,   9018, 0, 14,    375),
    @(111,    4313,    0,    3295,    952, 0,  1,     66),
    @(122,    2882,    0,    2658,    213, 0,  0,     11),
    @(145,    1308,    1,    1078,    198, 0,  0,     32)
)

foreach ($r in $plainRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
}

# --- 3. Rows whose country swaps with a neighbour, carrying new figures --
# row, Name, B, C, D, E, F, G, H
$swapRows = @(
    @(108, "Hungria",            4853, 40, 3590, 656, 0, 0, 607),
    @(109, "Guinea Ecuatorial",  4821,  0, 2182, 2556, 0, 0,  83),
    @(123, "Eslovaquia",         2801, 62, 1944,  826, 0, 0,  31),
    @(124, "Surinam",            2761,  0, 1830,  891, 0, 0,  40),
    @(213, "Montserrat",           13,  0,   12,    0, 0, 0,   1),
    @(214, "Islas Malvinas",       13,  0,   13,    0, 0, 0,   0)
)

foreach ($r in $swapRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]
    $ws.Cells.Item($row, 7).Value = $r[7]
    $ws.Cells.Item($row, 8).Value = $r[8]
}
